# Apply the two textual edits described by the diff:
#  1. Slide 1 (title slide): update the date run "2025.09.09." -> "2025.09.11."
#  2. Slide 32 ("A félév tematikája"): merge the three runs
#     "forradalma (" + "2020–22" + "): " into a single run "forradalma (2020–22): "

$p = $ppt.ActivePresentation

# --- Edit 1: slide 1, subtitle date ---------------------------------------
$slide1 = $p.Slides.Item(1)
$dateShape = $slide1.Shapes.Item(2)
$dateRange = $dateShape.TextFrame.TextRange
$dateRange.Paragraphs(1).Runs(1).Text = "2025.09.11."

# --- Edit 2: slide 32, merge "forradalma (2020-22): " runs ----------------
$slide32 = $p.Slides.Item(32)
$bodyShape = $slide32.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange

$enDash = [char]0x2013
$target = "forradalma (2020" + $enDash + "22): "

$fullText = $bodyRange.Text
$startIdx = $fullText.IndexOf($target)
if ($startIdx -ge 0) {
    $span = $bodyRange.Characters($startIdx + 1, $target.Length)
    $span.Text = $target
}
